# edit.ps1
# Applies the changes described by the diff:
#   1. Insert a collapsed "_GoBack" bookmark at the very start of the document
#      (right after the pPr of the first paragraph, before its first run).
#   2. Remove the "_GoBack" bookmark that currently splits the TOC heading run
#      ("Table of Co" / "ntents") and merge the text back into a single run
#      "Table of Contents".
#   3. Remove line-numbering (<w:lnNumType .../>) from the (only) section's
#      sectPr and add a <w:docGrid w:linePitch="326"/> element.
#   4. Update the cached page-number text in the second footer from "5" to "3".
#
# Strategy: pull the canonical WordOpenXML (flat-OPC) representation of the
# whole package, perform precise, uniqueness-checked string substitutions,
# then push the modified XML back via the WordOpenXML property. This avoids
# several quirks/limitations of the COM object model surfaced for this
# runtime (Bookmarks.Add mishandling collapsed ranges at offset 0,
# LineNumbering properties not round-tripping to XML, etc.).

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

function Assert-Count($haystack, $needle, $expected) {
    $count = 0
    $idx = 0
    while (($idx = $haystack.IndexOf($needle, $idx)) -ne -1) {
        $count++
        $idx += [Math]::Max($needle.Length, 1)
    }
    if ($count -ne $expected) {
        throw "Expected $expected occurrence(s) of [$needle] but found $count"
    }
}

# ---------------------------------------------------------------------
# Change 2: merge the TOC heading runs and drop the mid-word bookmark.
# ---------------------------------------------------------------------
$tocOld = '<w:r w:rsidRPr="00B72C39"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:sz w:val="32"/></w:rPr><w:t>Table of Co</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidRPr="00B72C39"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:sz w:val="32"/></w:rPr><w:t>ntents</w:t></w:r>'
$tocNew = '<w:r w:rsidRPr="00B72C39"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:sz w:val="32"/></w:rPr><w:t>Table of Contents</w:t></w:r>'
Assert-Count $xml $tocOld 1
$xml = $xml.Replace($tocOld, $tocNew)

# ---------------------------------------------------------------------
# Change 1: insert the "_GoBack" bookmark at the very start of the body,
# immediately after the first paragraph's pPr and before its first run.
# ---------------------------------------------------------------------
$introOld = '<w:ind w:left="0" w:right="0"/><w:jc w:val="right"/><w:rPr><w:i w:val="0"/><w:color w:val="4F81BD" w:themeColor="accent1"/><w:sz w:val="24"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:r><w:rPr><w:i w:val="0"/><w:color w:val="4F81BD" w:themeColor="accent1"/><w:sz w:val="24"/><w:szCs w:val="44"/></w:rPr><w:t>March 1, 2013</w:t></w:r>'
$introNew = '<w:ind w:left="0" w:right="0"/><w:jc w:val="right"/><w:rPr><w:i w:val="0"/><w:color w:val="4F81BD" w:themeColor="accent1"/><w:sz w:val="24"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:i w:val="0"/><w:color w:val="4F81BD" w:themeColor="accent1"/><w:sz w:val="24"/><w:szCs w:val="44"/></w:rPr><w:t>March 1, 2013</w:t></w:r>'
Assert-Count $xml $introOld 1
$xml = $xml.Replace($introOld, $introNew)

# ---------------------------------------------------------------------
# Change 3: drop line numbering from the section properties and add
# a docGrid element.
# ---------------------------------------------------------------------
$sectOld = '<w:lnNumType w:countBy="1" w:restart="continuous"/><w:cols w:space="720"/></w:sectPr>'
$sectNew = '<w:cols w:space="720"/><w:docGrid w:linePitch="326"/></w:sectPr>'
Assert-Count $xml $sectOld 1
$xml = $xml.Replace($sectOld, $sectNew)

# ---------------------------------------------------------------------
# Change 4: update the cached PAGE field text in the second footer
# (footer2.xml) from "5" to "3".
# ---------------------------------------------------------------------
$pageOld = '<w:rStyle w:val="PageNumber"/><w:noProof/></w:rPr><w:t>5</w:t>'
$pageNew = '<w:rStyle w:val="PageNumber"/><w:noProof/></w:rPr><w:t>3</w:t>'
Assert-Count $xml $pageOld 1
$xml = $xml.Replace($pageOld, $pageNew)

# ---------------------------------------------------------------------
# Push the modified XML back into the document.
# ---------------------------------------------------------------------
$d.WordOpenXML = $xml

Write-Host "edit.ps1 applied successfully"
